# 03_Catalogo_Productos - restructure product catalog columns
#  - Remove last two products (PR011, PR012)
#  - Replace/rename columns D..I and append new columns J, K, L
#  - Resize columns to the new widths

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the two products that are being discontinued (rows 12 and 13)
# ---------------------------------------------------------------------------
$ws.Rows("12:13").Delete()

# ---------------------------------------------------------------------------
# 2) Column widths (A..L)
#    COM ColumnWidth reads back ~0.83 lower than the stored OOXML <col width>
#    for this workbook's font, so subtract that offset to land exactly on
#    the target width after save.
# ---------------------------------------------------------------------------
$widths = @(8, 18, 25, 25, 30, 25, 12, 12, 12, 50, 20, 15)
for ($c = 1; $c -le $widths.Count; $c++) {
    $ws.Columns.Item($c).ColumnWidth = $widths[$c - 1] - 0.83
}

# ---------------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------------
$headers = @(
    "ID",
    "Nombre Producto",
    "Descripción",
    "Colores Asociados",
    "Flores Asociadas",
    "Tipos Macetero Posibles",
    "Vista (360/180)",
    "Tamaño (cm)",
    "Precio Venta",
    "Cuidados",
    "Foto",
    "Disponible Shopify"
)

for ($c = 1; $c -le $headers.Count; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    if ($c -gt 9) {
        # New header cells: clone the look of an existing header cell
        # (bold white font on blue fill) instead of re-building the format
        # by hand so the workbook reuses the same style record.
        $ws.Range("A1").Copy()
        $cell.PasteSpecial(-4122)
    }
    $cell.Value = $headers[$c - 1]
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Product rows (row 2 .. row 11)
#    Columns: A ID | B Nombre | C Descripcion | D Colores Asociados |
#             E Flores Asociadas | F Tipos Macetero Posibles |
#             G Vista 360/180 | H Tamaño (cm) | I Precio Venta |
#             J Cuidados | K Foto | L Disponible Shopify
# ---------------------------------------------------------------------------
$products = @(
    @("PR001","Pasión Roja","Arreglo elegante en tonos rojos","Rojo, Verde oscuro, Burdeo","Rosa roja, Clavel rojo, Eucalipto","Florero vidrio cilíndrico","360","25 x 35",35000,"Cambiar agua cada 2 días, evitar luz directa, cortar tallos en diagonal","passion-roja.jpg","Sí"),
    @("PR002","Sueño Blanco","Delicado arreglo en blancos puros","Blanco, Verde claro","Rosa blanca, Lirio blanco, Gerbera blanca","Florero vidrio redondo","360","22 x 30",32000,"Cambiar agua diariamente, mantener en lugar fresco, eliminar polen de lirios","sueno-blanco.jpg","Sí"),
    @("PR003","Jardín Primaveral","Mezcla de colores vibrantes","Amarillo, Naranja, Rosado, Morado","Gerbera, Alstroemeria, Rosa, Solidago","Florero vidrio grande","360","30 x 40",42000,"Cambiar agua cada 2-3 días, exposición luz indirecta, agregar nutriente floral","jardin-primaveral.jpg","Sí"),
    @("PR004","Elegancia Rosa","Rosas rosadas en florero","Rosado, Verde, Blanco","Rosa rosada, Eucalipto, Solidago","Florero cerámica blanco","180","20 x 35",38000,"Cambiar agua cada 2 días, cortar 1cm de tallo cada 3 días, ambiente fresco","elegancia-rosa.jpg","Sí"),
    @("PR005","Sol Radiante","Girasoles y flores amarillas","Amarillo, Naranja, Verde","Girasol, Gerbera naranja, Solidago","Florero vidrio cilíndrico","360","28 x 38",30000,"Cambiar agua diariamente, requiere buena luz, girasoles duran 7-10 días","sol-radiante.jpg","Sí"),
    @("PR006","Dulce Lirio","Lirios blancos y rosados","Blanco, Rosado, Verde","Lirio blanco, Lirio rosado, Eucalipto","Florero vidrio burbuja","360","25 x 40",45000,"Cambiar agua cada 2 días, quitar polen para evitar manchas, cortar en diagonal","dulce-lirio.jpg","Sí"),
    @("PR007","Campo Silvestre","Arreglo rústico en canasto","Multicolor natural","Mix de temporada, Gerbera, Alstroemeria, Follaje","Canasto mimbre rectangular","180","35 x 25",48000,"Verificar nivel de agua en esponja, rociar flores 1 vez al día","campo-silvestre.jpg","Sí"),
    @("PR008","Orquídea Imperial","Orquídeas blancas premium","Blanco puro","Orquídea phalaenopsis","Macetero cerámica gris","360","18 x 45",55000,"Regar 1 vez por semana, luz indirecta brillante, no mojar flores","orquidea-imperial.jpg","Sí"),
    @("PR009","Ramo Clásico","Ramo de rosas rojas","Rojo intenso, Verde","Rosa roja, Eucalipto","Sin contenedor (ramo)","360","Ø 25",28000,"Cortar tallos y colocar en agua inmediatamente, cambiar agua diariamente","ramo-clasico.jpg","Sí"),
    @("PR010","Amor Eterno","Caja con rosas","Rojo, Rosado suave","Rosa roja, Rosa rosada","Caja cuadrada","360","25 x 25",65000,"Rociar con agua 1-2 veces al día, no exponer al sol directo","amor-eterno.jpg","Sí")
)

for ($i = 0; $i -lt $products.Count; $i++) {
    $r = $i + 2
    $row = $products[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # G holds "360"/"180" as text (matches source inlineStr type), not a number
    $ws.Cells.Item($r, 7).Value = "'" + $row[6]
    $ws.Cells.Item($r, 7).Style = $ws.Cells.Item($r, 1).Style

    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
}
